# Applies the "added stuff at the bottom of the sheets" commit:
#  - fills in the pair_kind (J) column for the four practice rows
#  - appends a new "stim details" block describing audio/video/image needs

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Practice rows (2-5) were missing a pair_kind value of "generic" in column J
foreach ($r in 2..5) {
    $ws.Cells.Item($r, 10).Value = "generic"   # column J
}

# New section header
$ws.Range("A27").Value = "stim details"

# New header row describing the stim-detail columns
$ws.Range("A28").Value = "month"
$ws.Range("B28").Value = "word_type"
$ws.Range("C28").Value = "need_audio"
$ws.Range("D28").Value = "need_image"
$ws.Range("E28").Value = "word"
$ws.Range("F28").Value = "count"
$ws.Range("G28").Value = "find images"

# New data rows: month number (A) and word_type (B)
$stimRows = @(
    @{ Row = 29; Month = 6; WordType = "video" },
    @{ Row = 30; Month = 6; WordType = "video" },
    @{ Row = 31; Month = 7; WordType = "video" },
    @{ Row = 32; Month = 7; WordType = "video" },
    @{ Row = 33; Month = 6; WordType = "audio" },
    @{ Row = 34; Month = 6; WordType = "audio" },
    @{ Row = 35; Month = 7; WordType = "audio" },
    @{ Row = 36; Month = 7; WordType = "audio" }
)

foreach ($entry in $stimRows) {
    $ws.Cells.Item($entry.Row, 1).Value = $entry.Month
    $ws.Cells.Item($entry.Row, 2).Value = $entry.WordType
}
